# "now trying on event logic"
# - Move the worksheet selection from D8 to A12
# - Give row 12 an explicit (custom) height of 16.5pt
# - Give row 13 an explicit (custom) height of 40.5pt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select A12 (becomes the new <selection activeCell="A12" sqref="A12"/>)
$ws.Range("A12").Select()

# Explicit row heights (emit ht="..." customHeight="1" on save)
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 40.5
